$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Suite")

# Precondition text (TC1, TC2, TC3 share the same text) - B8, B19, B29
$ws.Range("B8").Value = "Administrador esta autenticado no sistema e tem permissao para alterar Gerente de Desempenho"
$ws.Range("B19").Value = "Administrador esta autenticado no sistema e tem permissao para alterar Gerente de Desempenho"
$ws.Range("B29").Value = "Administrador esta autenticado no sistema e tem permissao para alterar Gerente de Desempenho"

# Step 1 expected result - D10, D21, D31
$ws.Range("D10").Value = "SYSTEM exibe a listagem dos Perfis de Competencias cadastrados com a opcao 'Alterar Gerente' dentre as varias exibidas"
$ws.Range("D21").Value = "SYSTEM exibe a listagem dos Perfis de Competencias cadastrados com a opcao 'Alterar Gerente' dentre as varias exibidas"
$ws.Range("D31").Value = "SYSTEM exibe a listagem dos Perfis de Competencias cadastrados com a opcao 'Alterar Gerente' dentre as varias exibidas"

# Step 3 description - B12, B23, B33
$ws.Range("B12").Value = "Administrador preenche o campo 'Login do Novo Gerente de Desempenho' para o Perfil de Competencias"
$ws.Range("B23").Value = "Administrador preenche o campo 'Login do Novo Gerente de Desempenho' para o Perfil de Competencias"
$ws.Range("B33").Value = "Administrador preenche o campo 'Login do Novo Gerente de Desempenho' para o Perfil de Competencias"

# TC2 step 4 expected result (Cancelar) - D24
$ws.Range("D24").Value = "SYSTEM apresenta o Catalogo (Perfis) de Competencias sem nenhuma alteracao"
